$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.064.68"
$ws.Range("E2").Value = "  +0.30%  "
$ws.Range("D3").Value = "2.061.00"
$ws.Range("E3").Value = "  -0.96%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.37"
$ws.Range("E5").Value = "  -0.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.673"
$ws.Range("E6").Value = "  +2.41%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "55.83"
$ws.Range("E8").Value = "  +10.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "60.67"
$ws.Range("E9").Value = "  +0.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.379"
$ws.Range("E10").Value = "  +0.81%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0802"
$ws.Range("E11").Value = "  +7.80%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.107"
$ws.Range("E12").Value = "  +0.22%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.14"
$ws.Range("E13").Value = "  +3.74%  "
$ws.Range("D14").Value = "2.358.29"
$ws.Range("E14").Value = "  -1.25%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.813"
$ws.Range("E15").Value = "  -2.09%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.30"
$ws.Range("E16").Value = "  +2.81%  "
$ws.Range("D17").Value = "2.058.98"
$ws.Range("E17").Value = "  -0.76%  "
$ws.Range("D18").Value = "36.968.34"
$ws.Range("E18").Value = "  +0.29%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.28"
$ws.Range("E19").Value = "  +2.10%  "
$ws.Range("D20").Value = "0.0₃0919"
$ws.Range("E20").Value = "  +11.68%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.25"
$ws.Range("E21").Value = "  +7.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.38"
$ws.Range("E22").Value = "  +1.88%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.79"
$ws.Range("E23").Value = "  -0.99%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.42"
$ws.Range("E25").Value = "  -3.22%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "171.42"
$ws.Range("E26").Value = "  +1.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.10"
$ws.Range("E27").Value = "  -3.31%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.15"
$ws.Range("E28").Value = "  -3.53%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.00"
$ws.Range("E29").Value = "  +0.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.125"
$ws.Range("E30").Value = "  +1.94%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.59"
$ws.Range("E31").Value = "  +2.21%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.06"
$ws.Range("E32").Value = "  -8.44%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0627"
$ws.Range("E33").Value = "  +3.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.39"
$ws.Range("E34").Value = "  +7.27%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0866"
$ws.Range("E36").Value = "  -4.92%  "
$ws.Range("E37").Value = "  -2.75%  "
$ws.Range("E38").Value = "  -2.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.35"
$ws.Range("E39").Value = "  +2.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.105"
$ws.Range("E40").Value = "  +22.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.78"
$ws.Range("E41").Value = "  +65.27%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "18.18"
$ws.Range("E42").Value = "  +4.92%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0225"
$ws.Range("E43").Value = "  +0.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.15"
$ws.Range("E44").Value = "  -0.42%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "96.81"
$ws.Range("E45").Value = "  -0.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.78"
$ws.Range("E46").Value = "  -0.71%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.36"
$ws.Range("E47").Value = "  +12.15%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.44"
$ws.Range("E48").Value = "  +8.35%  "
$ws.Range("D49").Value = "1.297.91"
$ws.Range("E49").Value = "  -3.39%  "
$ws.Range("E50").Value = "  +0.36%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.84"
$ws.Range("E51").Value = "  -2.04%  "
